$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.551.82"
$ws.Range("E2").Value = "  +2.21%  "
$ws.Range("D3").Value = "1.790.29"
$ws.Range("E3").Value = "  +4.17%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'313.97"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.5349"
$ws.Range("E7").Value = "  +9.88%  "
$ws.Range("D8").Value = "'0.3762"
$ws.Range("E8").Value = "  +7.78%  "
$ws.Range("D9").Value = "'42.93"
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("D10").Value = "'0.07504"
$ws.Range("E10").Value = "  +3.46%  "
$ws.Range("D11").Value = "'1.113"
$ws.Range("E11").Value = "  +6.42%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "'20.87"
$ws.Range("E13").Value = "  +5.01%  "
$ws.Range("D14").Value = "'6.173"
$ws.Range("E14").Value = "  +5.46%  "
$ws.Range("D15").Value = "1.788.79"
$ws.Range("E15").Value = "  +3.77%  "
$ws.Range("D16").Value = "'7.080"
$ws.Range("D17").Value = "'90.82"
$ws.Range("D18").Value = "'0.00001068"
$ws.Range("E18").Value = "  +3.03%  "
$ws.Range("D19").Value = "'0.06500"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").Value = "'0.9996"
$ws.Range("D21").Value = "'16.94"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("D22").Value = "'5.933"
$ws.Range("E22").Value = "  +5.28%  "
$ws.Range("D23").Value = "27.591.98"
$ws.Range("E23").Value = "  +2.16%  "
$ws.Range("E24").Value = "  +3.93%  "
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("D26").Value = "'20.47"
$ws.Range("E26").Value = "  +2.76%  "
$ws.Range("D27").Value = "'155.30"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("D28").Value = "'2.384"
$ws.Range("E28").Value = "  +15.50%  "
$ws.Range("D29").Value = "1.993.29"
$ws.Range("E29").Value = "  +3.84%  "
$ws.Range("D30").Value = "'121.84"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").Value = "'1.118"
$ws.Range("E31").Value = "  +8.83%  "
$ws.Range("D32").Value = "'0.1029"
$ws.Range("E32").Value = "  +10.97%  "
$ws.Range("D33").Value = "'5.674"
$ws.Range("E33").Value = "  +6.13%  "
$ws.Range("D34").Value = "'3.610"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("D35").Value = "'0.02279"
$ws.Range("E35").Value = "  +4.72%  "
$ws.Range("D36").Value = "'8.671"
$ws.Range("E36").Value = "  +15.05%  "
$ws.Range("D37").Value = "'0.06020"
$ws.Range("E37").Value = "  +2.18%  "
$ws.Range("D38").Value = "'4.987"
$ws.Range("E38").Value = "  +5.26%  "
$ws.Range("D39").Value = "'0.2082"
$ws.Range("E39").Value = "  +4.47%  "
$ws.Range("E40").Value = "  +3.96%  "
$ws.Range("D41").Value = "'0.6244"
$ws.Range("E41").Value = "  +4.34%  "
$ws.Range("D42").Value = "'1.413"
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("D43").Value = "'0.9995"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "'1.143"
$ws.Range("E44").Value = "  +4.92%  "
$ws.Range("D45").Value = "'13.36"
$ws.Range("E45").Value = "  +4.47%  "
$ws.Range("D46").Value = "'0.5859"
$ws.Range("E46").Value = "  +4.29%  "
$ws.Range("D47").Value = "'3.635"
$ws.Range("D48").Value = "'121.45"
$ws.Range("E48").Value = "  +3.19%  "
$ws.Range("D49").Value = "'1.914"
$ws.Range("E49").Value = "  +4.46%  "
$ws.Range("D50").Value = "'1.132"
$ws.Range("E50").Value = "  +2.33%  "
$ws.Range("E51").Value = "  +1.63%  "
